# Update countries & provincias Spain
#
# Refreshes the COVID-19 snapshot: the "last updated" timestamp, a handful
# of per-country totals (rows for Estados Unidos/USA, Reino Unido/UK,
# Francia/France, India - rows 4, 8, 11, 14), and the block of rows 69-72
# where "Azerbaiyan" overtook "Camerun" in the case-count sort order
# (so Azerbaiyan/Camerun/Irak/Hungria each shift down one slot while
# picking up their refreshed figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value2 = 'Datos actualizados a 21 de Mayo de 2020 a las 15:35'

# Estados Unidos (row 4)
$ws.Range("B4").Value2 = 1595081
$ws.Range("C4").Value2 = 2358
$ws.Range("E4").Value2 = 1129094
$ws.Range("G4").Value2 = 80
$ws.Range("H4").Value2 = 95016

# Reino Unido (row 8)
$ws.Range("B8").Value2 = 250908
$ws.Range("C8").Value2 = 2615
$ws.Range("G8").Value2 = 338
$ws.Range("H8").Value2 = 36042

# Francia (row 11)
$ws.Range("B11").Value2 = 178568
$ws.Range("C11").Value2 = 37
$ws.Range("E11").Value2 = 12297
$ws.Range("G11").Value2 = 1
$ws.Range("H11").Value2 = 8271

# India (row 14)
$ws.Range("B14").Value2 = 113461
$ws.Range("C14").Value2 = 1433
$ws.Range("D14").Value2 = 46002
$ws.Range("E14").Value2 = 64002
$ws.Range("G14").Value2 = 23
$ws.Range("H14").Value2 = 3457

# Azerbaiyan moves ahead of Camerun; Camerun/Irak/Hungria each shift down
# one row and bring their (unchanged) figures with them - Sudan (row 73)
# is unaffected.
$ws.Range("A69").Value2 = 'Azerbaiyan'
$ws.Range("B69").Value2 = 3749
$ws.Range("C69").Value2 = 118
$ws.Range("D69").Value2 = 2340
$ws.Range("E69").Value2 = 1365
$ws.Range("G69").Value2 = 1
$ws.Range("H69").Value2 = 44

$ws.Range("A70").Value2 = 'Camerun'
$ws.Range("B70").Value2 = 3733
$ws.Range("D70").Value2 = 1595
$ws.Range("E70").Value2 = 1992
$ws.Range("H70").Value2 = 146

$ws.Range("A71").Value2 = 'Irak'
$ws.Range("B71").Value2 = 3724
$ws.Range("C71").Value2 = 0
$ws.Range("D71").Value2 = 2438
$ws.Range("E71").Value2 = 1152
$ws.Range("G71").Value2 = 0
$ws.Range("H71").Value2 = 134

$ws.Range("A72").Value2 = 'Hungria'
$ws.Range("B72").Value2 = 3641
$ws.Range("C72").Value2 = 43
$ws.Range("D72").Value2 = 1509
$ws.Range("E72").Value2 = 1659
$ws.Range("G72").Value2 = 3
$ws.Range("H72").Value2 = 473
